$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "hepatitis A" / "Hepatitis B" symptom-combination rows (145:156).
# They are superseded by the new Stroke / Coronary Artery Disease / IBD / Migraine blocks below.
$ws.Range("A145:Q156").ClearContents()

# Add the new sentence models (disease/symptom combination rows) for:
#   Stroke, Coronary Artery Disease, Inflammatory Bowel Disease (IBD), Migraine
# row 132
$ws.Range("A132").Value = 'Stroke'
$ws.Range("B132").Value = 'vision_impairment'
$ws.Range("C132").Value = 'headache'
$ws.Range("D132").Value = 'numbness'
$ws.Range("E132").Value = 'confusion  '
$ws.Range("F132").Value = 'trouble_walking'
$ws.Range("G132").Value = 'trouble_speaking'
$ws.Range("H132").Value = 'facial_paralysis'
$ws.Range("I132").Value = 'feeling_dizzy'

# row 133
$ws.Range("A133").Value = 'Stroke'
$ws.Range("B133").Value = 'vision_impairment'
$ws.Range("C133").Value = 'headache'
$ws.Range("D133").Value = 'numbness'
$ws.Range("E133").Value = 'confusion  '
$ws.Range("F133").Value = 'trouble_walking'
$ws.Range("G133").Value = 'trouble_speaking'
$ws.Range("H133").Value = 'facial_paralysis'
$ws.Range("I133").Value = 'falling_over'

# row 134
$ws.Range("A134").Value = 'Stroke'
$ws.Range("B134").Value = 'vision_impairment'
$ws.Range("C134").Value = 'headache'
$ws.Range("D134").Value = 'numbness'
$ws.Range("E134").Value = 'confusion  '
$ws.Range("F134").Value = 'trouble_walking'
$ws.Range("G134").Value = 'trouble_speaking'
$ws.Range("H134").Value = 'arm_paralysis'
$ws.Range("I134").Value = 'feeling_dizzy'

# row 135
$ws.Range("A135").Value = 'Stroke'
$ws.Range("B135").Value = 'vision_impairment'
$ws.Range("C135").Value = 'headache'
$ws.Range("D135").Value = 'numbness'
$ws.Range("E135").Value = 'confusion  '
$ws.Range("F135").Value = 'trouble_walking'
$ws.Range("G135").Value = 'trouble_speaking'
$ws.Range("H135").Value = 'arm_paralysis'
$ws.Range("I135").Value = 'falling_over'

# row 136
$ws.Range("A136").Value = 'Stroke'
$ws.Range("B136").Value = 'vision_impairment'
$ws.Range("C136").Value = 'headache'
$ws.Range("D136").Value = 'numbness'
$ws.Range("E136").Value = 'confusion  '
$ws.Range("F136").Value = 'trouble_walking'
$ws.Range("G136").Value = 'trouble_speaking'
$ws.Range("H136").Value = 'leg_paralysis'
$ws.Range("I136").Value = 'feeling_dizzy'

# row 137
$ws.Range("A137").Value = 'Stroke'
$ws.Range("B137").Value = 'vision_impairment'
$ws.Range("C137").Value = 'headache'
$ws.Range("D137").Value = 'numbness'
$ws.Range("E137").Value = 'confusion  '
$ws.Range("F137").Value = 'trouble_walking'
$ws.Range("G137").Value = 'trouble_speaking'
$ws.Range("H137").Value = 'leg_paralysis'
$ws.Range("I137").Value = 'falling_over'

# row 138
$ws.Range("A138").Value = 'Stroke'
$ws.Range("B138").Value = 'vision_impairment'
$ws.Range("C138").Value = 'headache'
$ws.Range("D138").Value = 'numbness'
$ws.Range("E138").Value = 'memory_loss'
$ws.Range("F138").Value = 'trouble_walking'
$ws.Range("G138").Value = 'trouble_speaking'
$ws.Range("H138").Value = 'facial_paralysis'
$ws.Range("I138").Value = 'feeling_dizzy'

# row 139
$ws.Range("A139").Value = 'Stroke'
$ws.Range("B139").Value = 'vision_impairment'
$ws.Range("C139").Value = 'headache'
$ws.Range("D139").Value = 'numbness'
$ws.Range("E139").Value = 'memory_loss'
$ws.Range("F139").Value = 'trouble_walking'
$ws.Range("G139").Value = 'trouble_speaking'
$ws.Range("H139").Value = 'facial_paralysis'
$ws.Range("I139").Value = 'falling_over'

# row 140
$ws.Range("A140").Value = 'Stroke'
$ws.Range("B140").Value = 'vision_impairment'
$ws.Range("C140").Value = 'headache'
$ws.Range("D140").Value = 'numbness'
$ws.Range("E140").Value = 'memory_loss'
$ws.Range("F140").Value = 'trouble_walking'
$ws.Range("G140").Value = 'trouble_speaking'
$ws.Range("H140").Value = 'arm_paralysis'
$ws.Range("I140").Value = 'feeling_dizzy'

# row 141
$ws.Range("A141").Value = 'Stroke'
$ws.Range("B141").Value = 'vision_impairment'
$ws.Range("C141").Value = 'headache'
$ws.Range("D141").Value = 'numbness'
$ws.Range("E141").Value = 'memory_loss'
$ws.Range("F141").Value = 'trouble_walking'
$ws.Range("G141").Value = 'trouble_speaking'
$ws.Range("H141").Value = 'arm_paralysis'
$ws.Range("I141").Value = 'falling_over'

# row 142
$ws.Range("A142").Value = 'Stroke'
$ws.Range("B142").Value = 'vision_impairment'
$ws.Range("C142").Value = 'headache'
$ws.Range("D142").Value = 'numbness'
$ws.Range("E142").Value = 'memory_loss'
$ws.Range("F142").Value = 'trouble_walking'
$ws.Range("G142").Value = 'trouble_speaking'
$ws.Range("H142").Value = 'leg_paralysis'
$ws.Range("I142").Value = 'feeling_dizzy'

# row 143
$ws.Range("A143").Value = 'Stroke'
$ws.Range("B143").Value = 'vision_impairment'
$ws.Range("C143").Value = 'headache'
$ws.Range("D143").Value = 'numbness'
$ws.Range("E143").Value = 'memory_loss'
$ws.Range("F143").Value = 'trouble_walking'
$ws.Range("G143").Value = 'trouble_speaking'
$ws.Range("H143").Value = 'leg_paralysis'
$ws.Range("I143").Value = 'falling_over'

# row 144
$ws.Range("A144").Value = 'Coronary Artery Disease'
$ws.Range("B144").Value = 'angina'
$ws.Range("C144").Value = 'shortness_of_breath'
$ws.Range("E144").Value = 'feeling_faint'
$ws.Range("F144").Value = 'nausea'
$ws.Range("G144").Value = 'pain_in_your_neck'

# row 145
$ws.Range("A145").Value = 'Coronary Artery Disease'
$ws.Range("B145").Value = 'angina'
$ws.Range("C145").Value = 'shortness_of_breath'
$ws.Range("E145").Value = 'feeling_faint'
$ws.Range("F145").Value = 'nausea'
$ws.Range("G145").Value = 'pain_in_your_shoulders'

# row 146
$ws.Range("A146").Value = 'Coronary Artery Disease'
$ws.Range("B146").Value = 'angina'
$ws.Range("C146").Value = 'shortness_of_breath'
$ws.Range("E146").Value = 'feeling_faint'
$ws.Range("F146").Value = 'nausea'
$ws.Range("G146").Value = 'pain_in_your_neck'
$ws.Range("H146").Value = 'pain_in_your_jaw'

# row 147
$ws.Range("A147").Value = 'Coronary Artery Disease'
$ws.Range("B147").Value = 'angina'
$ws.Range("C147").Value = 'shortness_of_breath'
$ws.Range("E147").Value = 'feeling_faint'
$ws.Range("F147").Value = 'nausea'
$ws.Range("G147").Value = 'pain_in_your_shoulders'
$ws.Range("H147").Value = 'pain_in_your_jaw'

# row 148
$ws.Range("A148").Value = 'Coronary Artery Disease'
$ws.Range("B148").Value = 'angina'
$ws.Range("C148").Value = 'shortness_of_breath'
$ws.Range("E148").Value = 'feeling_faint'
$ws.Range("F148").Value = 'nausea'
$ws.Range("G148").Value = 'pain_in_your_jaw'

# row 149
$ws.Range("A149").Value = 'Coronary Artery Disease'
$ws.Range("B149").Value = 'angina'
$ws.Range("C149").Value = 'shortness_of_breath'
$ws.Range("E149").Value = 'feeling_faint'
$ws.Range("F149").Value = 'nausea'
$ws.Range("G149").Value = 'pain_in_your_arms'

# row 150
$ws.Range("A150").Value = 'Coronary Artery Disease'
$ws.Range("B150").Value = 'angina'
$ws.Range("C150").Value = 'shortness_of_breath'
$ws.Range("E150").Value = 'feeling_faint'
$ws.Range("F150").Value = 'nausea'
$ws.Range("G150").Value = 'pain_in_your_neck'
$ws.Range("H150").Value = 'pain_in_your_shoulders'
$ws.Range("I150").Value = 'pain_in_your_jaw'

# row 151
$ws.Range("A151").Value = 'Coronary Artery Disease'
$ws.Range("B151").Value = 'angina'
$ws.Range("C151").Value = 'shortness_of_breath'
$ws.Range("E151").Value = 'feeling_faint'
$ws.Range("F151").Value = 'nausea'
$ws.Range("G151").Value = 'pain_in_your_neck'
$ws.Range("H151").Value = 'pain_in_your_shoulders'
$ws.Range("I151").Value = 'pain_in_your_arms'

# row 152
$ws.Range("A152").Value = 'Coronary Artery Disease'
$ws.Range("B152").Value = 'angina'
$ws.Range("C152").Value = 'shortness_of_breath'
$ws.Range("E152").Value = 'feeling_faint'
$ws.Range("F152").Value = 'cold_sweat'
$ws.Range("G152").Value = 'pain_in_your_neck'

# row 153
$ws.Range("A153").Value = 'Coronary Artery Disease'
$ws.Range("B153").Value = 'angina'
$ws.Range("C153").Value = 'shortness_of_breath'
$ws.Range("E153").Value = 'feeling_faint'
$ws.Range("F153").Value = 'cold_sweat'
$ws.Range("G153").Value = 'pain_in_your_shoulders'

# row 154
$ws.Range("A154").Value = 'Coronary Artery Disease'
$ws.Range("B154").Value = 'angina'
$ws.Range("C154").Value = 'shortness_of_breath'
$ws.Range("E154").Value = 'feeling_faint'
$ws.Range("F154").Value = 'cold_sweat'
$ws.Range("G154").Value = 'pain_in_your_neck'
$ws.Range("H154").Value = 'pain_in_your_jaw'

# row 155
$ws.Range("A155").Value = 'Coronary Artery Disease'
$ws.Range("B155").Value = 'angina'
$ws.Range("C155").Value = 'shortness_of_breath'
$ws.Range("E155").Value = 'feeling_faint'
$ws.Range("F155").Value = 'cold_sweat'
$ws.Range("G155").Value = 'pain_in_your_shoulders'
$ws.Range("H155").Value = 'pain_in_your_jaw'

# row 156
$ws.Range("A156").Value = 'Coronary Artery Disease'
$ws.Range("B156").Value = 'angina'
$ws.Range("C156").Value = 'shortness_of_breath'
$ws.Range("E156").Value = 'feeling_faint'
$ws.Range("F156").Value = 'cold_sweat'
$ws.Range("G156").Value = 'pain_in_your_jaw'

# row 157
$ws.Range("A157").Value = 'Coronary Artery Disease'
$ws.Range("B157").Value = 'angina'
$ws.Range("C157").Value = 'shortness_of_breath'
$ws.Range("E157").Value = 'feeling_faint'
$ws.Range("F157").Value = 'cold_sweat'
$ws.Range("G157").Value = 'pain_in_your_arms'

# row 158
$ws.Range("A158").Value = 'Coronary Artery Disease'
$ws.Range("B158").Value = 'angina'
$ws.Range("C158").Value = 'shortness_of_breath'
$ws.Range("E158").Value = 'feeling_faint'
$ws.Range("F158").Value = 'cold_sweat'
$ws.Range("G158").Value = 'pain_in_your_neck'
$ws.Range("H158").Value = 'pain_in_your_shoulders'
$ws.Range("I158").Value = 'pain_in_your_jaw'

# row 159
$ws.Range("A159").Value = 'Coronary Artery Disease'
$ws.Range("B159").Value = 'angina'
$ws.Range("C159").Value = 'shortness_of_breath'
$ws.Range("E159").Value = 'feeling_faint'
$ws.Range("F159").Value = 'cold_sweat'
$ws.Range("G159").Value = 'pain_in_your_neck'
$ws.Range("H159").Value = 'pain_in_your_shoulders'
$ws.Range("I159").Value = 'pain_in_your_arms'

# row 160
$ws.Range("A160").Value = 'Inflammatory Bowel Disease (IBD)'
$ws.Range("B160").Value = 'tummy_pain'
$ws.Range("C160").Value = 'diarrhea'
$ws.Range("D160").Value = 'loss_of_appetite'
$ws.Range("E160").Value = 'weight_loss'
$ws.Range("F160").Value = 'extreme_tiredness'
$ws.Range("G160").Value = 'blood_in_stool'

# row 161
$ws.Range("A161").Value = 'Inflammatory Bowel Disease (IBD)'
$ws.Range("B161").Value = 'tummy_pain'
$ws.Range("C161").Value = 'diarrhea'
$ws.Range("D161").Value = 'loss_of_appetite'
$ws.Range("E161").Value = 'weight_loss'
$ws.Range("F161").Value = 'extreme_tiredness'
$ws.Range("G161").Value = 'mucus_in_stool'

# row 162
$ws.Range("A162").Value = 'Inflammatory Bowel Disease (IBD)'
$ws.Range("B162").Value = 'tummy_pain'
$ws.Range("C162").Value = 'diarrhea'
$ws.Range("D162").Value = 'loss_of_appetite'
$ws.Range("E162").Value = 'weight_loss'
$ws.Range("F162").Value = 'extreme_tiredness'
$ws.Range("G162").Value = 'blood_in_stool'
$ws.Range("H162").Value = 'joint_pain'

# row 163
$ws.Range("A163").Value = 'Inflammatory Bowel Disease (IBD)'
$ws.Range("B163").Value = 'tummy_pain'
$ws.Range("C163").Value = 'diarrhea'
$ws.Range("D163").Value = 'loss_of_appetite'
$ws.Range("E163").Value = 'weight_loss'
$ws.Range("F163").Value = 'extreme_tiredness'
$ws.Range("G163").Value = 'mucus_in_stool'
$ws.Range("H163").Value = 'night_sweats'

# row 164
$ws.Range("A164").Value = 'Inflammatory Bowel Disease (IBD)'
$ws.Range("B164").Value = 'tummy_pain'
$ws.Range("C164").Value = 'diarrhea'
$ws.Range("D164").Value = 'loss_of_appetite'
$ws.Range("E164").Value = 'weight_loss'
$ws.Range("F164").Value = 'extreme_tiredness'
$ws.Range("G164").Value = 'blood_in_stool'
$ws.Range("H164").Value = 'skin_rash'

# row 165
$ws.Range("A165").Value = 'Inflammatory Bowel Disease (IBD)'
$ws.Range("B165").Value = 'tummy_pain'
$ws.Range("C165").Value = 'diarrhea'
$ws.Range("D165").Value = 'loss_of_appetite'
$ws.Range("E165").Value = 'weight_loss'
$ws.Range("F165").Value = 'extreme_tiredness'
$ws.Range("G165").Value = 'mucus_in_stool'
$ws.Range("H165").Value = 'fever'

# row 166
$ws.Range("A166").Value = 'Inflammatory Bowel Disease (IBD)'
$ws.Range("B166").Value = 'tummy_pain'
$ws.Range("C166").Value = 'diarrhea'
$ws.Range("D166").Value = 'loss_of_appetite'
$ws.Range("E166").Value = 'weight_loss'
$ws.Range("F166").Value = 'extreme_tiredness'
$ws.Range("G166").Value = 'blood_in_stool'
$ws.Range("H166").Value = 'joint_pain'
$ws.Range("I166").Value = 'night_sweats'

# row 167
$ws.Range("A167").Value = 'Inflammatory Bowel Disease (IBD)'
$ws.Range("B167").Value = 'tummy_pain'
$ws.Range("C167").Value = 'diarrhea'
$ws.Range("D167").Value = 'loss_of_appetite'
$ws.Range("E167").Value = 'weight_loss'
$ws.Range("F167").Value = 'extreme_tiredness'
$ws.Range("G167").Value = 'mucus_in_stool'
$ws.Range("H167").Value = 'skin_rash'
$ws.Range("I167").Value = 'fever'

# row 168
$ws.Range("A168").Value = 'Inflammatory Bowel Disease (IBD)'
$ws.Range("B168").Value = 'tummy_pain'
$ws.Range("C168").Value = 'diarrhea'
$ws.Range("D168").Value = 'loss_of_appetite'
$ws.Range("E168").Value = 'weight_loss'
$ws.Range("F168").Value = 'extreme_tiredness'
$ws.Range("G168").Value = 'blood_in_stool'
$ws.Range("H168").Value = 'joint_pain'
$ws.Range("I168").Value = 'night_sweats'
$ws.Range("J168").Value = 'skin_rash'
$ws.Range("K168").Value = 'fever'

# row 169
$ws.Range("A169").Value = 'Inflammatory Bowel Disease (IBD)'
$ws.Range("B169").Value = 'tummy_pain'
$ws.Range("C169").Value = 'diarrhea'
$ws.Range("D169").Value = 'loss_of_appetite'
$ws.Range("E169").Value = 'weight_loss'
$ws.Range("F169").Value = 'extreme_tiredness'
$ws.Range("G169").Value = 'mucus_in_stool'
$ws.Range("H169").Value = 'joint_pain'
$ws.Range("I169").Value = 'night_sweats'
$ws.Range("J169").Value = 'skin_rash'
$ws.Range("K169").Value = 'fever'

# row 170
$ws.Range("A170").Value = 'Migraine'
$ws.Range("B170").Value = 'stiff_neck'
$ws.Range("C170").Value = 'excessive_hunger'
$ws.Range("D170").Value = 'visual_disturbances'
$ws.Range("E170").Value = 'dizziness'
$ws.Range("F170").Value = 'indigestion'
$ws.Range("G170").Value = 'headache'
$ws.Range("H170").Value = 'blurred_and_distorted_vision'

# row 171
$ws.Range("A171").Value = 'Migraine'
$ws.Range("B171").Value = 'stiff_neck'
$ws.Range("C171").Value = 'excessive_hunger'
$ws.Range("D171").Value = 'visual_disturbances'
$ws.Range("E171").Value = 'dizziness'
$ws.Range("F171").Value = 'indigestion'
$ws.Range("G171").Value = 'headache'
$ws.Range("H171").Value = 'blurred_and_distorted_vision'
$ws.Range("I171").Value = 'changes_in_the_mood'

# row 172
$ws.Range("A172").Value = 'Migraine'
$ws.Range("B172").Value = 'stiff_neck'
$ws.Range("C172").Value = 'excessive_hunger'
$ws.Range("D172").Value = 'visual_disturbances'
$ws.Range("E172").Value = 'dizziness'
$ws.Range("F172").Value = 'indigestion'
$ws.Range("G172").Value = 'headache'
$ws.Range("H172").Value = 'blurred_and_distorted_vision'
$ws.Range("I172").Value = 'difficulty_speaking'

# row 173
$ws.Range("A173").Value = 'Migraine'
$ws.Range("B173").Value = 'stiff_neck'
$ws.Range("C173").Value = 'excessive_hunger'
$ws.Range("D173").Value = 'visual_disturbances'
$ws.Range("E173").Value = 'dizziness'
$ws.Range("F173").Value = 'indigestion'
$ws.Range("G173").Value = 'headache'
$ws.Range("H173").Value = 'blurred_and_distorted_vision'
$ws.Range("I173").Value = 'numbness'

# row 174
$ws.Range("A174").Value = 'Migraine'
$ws.Range("B174").Value = 'stiff_neck'
$ws.Range("C174").Value = 'excessive_hunger'
$ws.Range("D174").Value = 'visual_disturbances'
$ws.Range("E174").Value = 'dizziness'
$ws.Range("F174").Value = 'indigestion'
$ws.Range("G174").Value = 'headache'
$ws.Range("H174").Value = 'blurred_and_distorted_vision'
$ws.Range("I174").Value = 'changes_in_the_mood'
$ws.Range("J174").Value = 'difficulty_speaking'

# row 175
$ws.Range("A175").Value = 'Migraine'
$ws.Range("B175").Value = 'stiff_neck'
$ws.Range("C175").Value = 'excessive_hunger'
$ws.Range("D175").Value = 'visual_disturbances'
$ws.Range("E175").Value = 'dizziness'
$ws.Range("F175").Value = 'indigestion'
$ws.Range("G175").Value = 'headache'
$ws.Range("H175").Value = 'blurred_and_distorted_vision'
$ws.Range("I175").Value = 'changes_in_the_mood'
$ws.Range("J175").Value = 'numbness'

# row 176
$ws.Range("A176").Value = 'Migraine'
$ws.Range("B176").Value = 'stiff_neck'
$ws.Range("C176").Value = 'excessive_hunger'
$ws.Range("D176").Value = 'visual_disturbances'
$ws.Range("E176").Value = 'dizziness'
$ws.Range("F176").Value = 'indigestion'
$ws.Range("G176").Value = 'headache'
$ws.Range("H176").Value = 'blurred_and_distorted_vision'
$ws.Range("I176").Value = 'changes_in_the_mood'
$ws.Range("J176").Value = 'difficulty_speaking'
$ws.Range("K176").Value = 'numbness'

